$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.546.02"
$ws.Range("E2").Value = "  +1.36%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.173.69"
$ws.Range("E3").Value = "  -0.93%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.53"
$ws.Range("E5").Value = "  -0.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.63"
$ws.Range("E6").Value = "  -0.51%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.170.64"
$ws.Range("E8").Value = "  -0.94%  "

$ws.Range("E9").Value = "  +1.64%  "

$ws.Range("E10").Value = "  -1.30%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.34"
$ws.Range("E11").Value = "  -0.40%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.456"
$ws.Range("E12").Value = "  +0.01%  "

$ws.Range("E13").Value = "  +0.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.71"
$ws.Range("E14").Value = "  +3.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.697.21"
$ws.Range("E15").Value = "  -0.91%  "

$ws.Range("E16").Value = "  -0.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.176.57"
$ws.Range("E17").Value = "  -0.77%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.585.58"
$ws.Range("E18").Value = "  +1.24%  "

$ws.Range("E19").Value = "  -2.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "462.67"
$ws.Range("E20").Value = "  -0.61%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.98"
$ws.Range("E21").Value = "  -0.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.697"
$ws.Range("E22").Value = "  -2.59%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.68"
$ws.Range("E23").Value = "  -0.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.25"
$ws.Range("E24").Value = "  -2.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.07"
$ws.Range("E25").Value = "  -0.75%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.16%  "

$ws.Range("E27").Value = "  -1.97%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.07"
$ws.Range("E29").Value = "  -0.63%  "

$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.74"
$ws.Range("E30").Value = "  -2.80%  "

$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.81"
$ws.Range("E31").Value = "  -1.96%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.33"
$ws.Range("E32").Value = "  -1.16%  "

$ws.Range("E33").Value = "  -1.69%  "

$ws.Range("E34").Value = "  -1.85%  "

$ws.Range("E35").Value = "  -2.36%  "

$ws.Range("E36").Value = "  +0.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.48"
$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0733"
$ws.Range("E38").Value = "  +4.26%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0391"
$ws.Range("E39").Value = "  -1.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.13"
$ws.Range("E40").Value = "  -0.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.113"
$ws.Range("E41").Value = "  -2.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.65"
$ws.Range("E42").Value = "  -0.30%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "393.37"
$ws.Range("E43").Value = "  -6.53%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.795.63"
$ws.Range("E44").Value = "  -7.39%  "

$ws.Range("E45").Value = "  -1.23%  "

$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "127.93"
$ws.Range("E46").Value = "  +1.75%  "

$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.97"
$ws.Range("E47").Value = "  -1.21%  "

$ws.Range("E48").Value = "  +0.02%  "

$ws.Range("E49").Value = "  -2.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.26"
$ws.Range("E50").Value = "  -3.07%  "

$ws.Range("E51").Value = "  -0.84%  "
